$d = $word.ActiveDocument

# Update the title/date paragraph
$d.Content.Find.Execute("2024-11-29 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-11-30 Saturday", 2)

# Update the division problems in the table (in document order)
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "18÷9=2, 0"
$t.Cell(1,2).Range.Text = "39÷7=5, 4"
$t.Cell(1,3).Range.Text = "87÷7=12, 3"
$t.Cell(1,4).Range.Text = "57÷3=19, 0"
$t.Cell(1,5).Range.Text = "87÷6=14, 3"
$t.Cell(5,1).Range.Text = "73÷2=36, 1"
$t.Cell(5,2).Range.Text = "58÷9=6, 4"
$t.Cell(5,3).Range.Text = "80÷7=11, 3"
$t.Cell(5,4).Range.Text = "80÷8=10, 0"
$t.Cell(5,5).Range.Text = "39÷9=4, 3"
$t.Cell(9,1).Range.Text = "34÷6=5, 4"
$t.Cell(9,2).Range.Text = "93÷6=15, 3"
$t.Cell(9,3).Range.Text = "39÷7=5, 4"
$t.Cell(9,4).Range.Text = "44÷2=22, 0"
$t.Cell(9,5).Range.Text = "88÷9=9, 7"
$t.Cell(13,1).Range.Text = "83÷6=13, 5"
$t.Cell(13,2).Range.Text = "57÷2=28, 1"
$t.Cell(13,3).Range.Text = "44÷9=4, 8"
$t.Cell(13,4).Range.Text = "30÷9=3, 3"
$t.Cell(13,5).Range.Text = "39÷4=9, 3"
$t.Cell(17,1).Range.Text = "31÷3=10, 1"
$t.Cell(17,2).Range.Text = "17÷2=8, 1"
$t.Cell(17,3).Range.Text = "60÷8=7, 4"
$t.Cell(17,4).Range.Text = "90÷2=45, 0"
$t.Cell(17,5).Range.Text = "43÷7=6, 1"
